# number-showcase.xlsx update:
#   - "#system" sheet, column V ("web" picklist) gains a new alphabetically
#     sorted entry "clickWithKeys(locator,keys)" at row 51, pushing the
#     previous rows 51..117 down to 52..118.
#   - "#system" sheet, cell C2 ("aws.ses" picklist) now reads
#     "sendHtmlMail(profile,to,subject,body)" instead of "sendMail(...)".
#   - the "web" defined name is widened from $V$2:$V$117 to $V$2:$V$118 to
#     cover the extra row.
#
# NOTE: a plain Range.Insert(xlShiftDown) on a single cell shifts the whole
# physical row (every column), which would incorrectly drag column F (the
# unrelated "desktop" picklist) down too. Column F must stay put, so the V
# column values are shifted manually, one cell at a time, instead of using
# Insert/EntireRow.Insert.

$wb = $excel.ActiveWorkbook
$sys = $wb.Worksheets.Item("#system")

$webCol = 22   # column V
$firstDataRow = 2
$oldLastRow = 117
$newLastRow = 118
$insertRow = 51

# Shift V(insertRow..oldLastRow) down into V(insertRow+1..newLastRow),
# walking from the bottom up so values aren't clobbered before they're read.
for ($r = $oldLastRow; $r -ge $insertRow; $r--) {
    $val = $sys.Cells.Item($r, $webCol).Value2
    $sys.Cells.Item($r + 1, $webCol).Value2 = $val
}

# New alphabetically-ordered entry: "clickByLabelAndWait(...)" (50) <
# "clickWithKeys(...)" < "close()" (now row 52).
$sys.Cells.Item($insertRow, $webCol).Value2 = "clickWithKeys(locator,keys)"

# aws.ses picklist: replace sendMail(...) with the new sendHtmlMail(...).
$sys.Cells.Item(2, 3).Value2 = "sendHtmlMail(profile,to,subject,body)"

# Widen the "web" named range to include the newly added row.
$webName = $wb.Names.Item("web")
$webName.RefersTo = "='#system'!`$V`$$firstDataRow`:`$V`$$newLastRow"
